$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 12) of portfolio data, matching the text-based
# date format already used in column A of the existing rows. Forcing the
# cell to Text before assignment keeps "2025-08-27" as a literal string
# instead of Excel auto-converting it to a date serial; ClearFormats()
# afterwards drops the temporary Text number format again so the cell
# ends up with the same (default) style as the rest of the data rows.
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2025-08-27"
$ws.Range("A12").ClearFormats()

$ws.Range("B12").Value = 57.09999847412109
$ws.Range("C12").Value = 680.5499877929688
$ws.Range("D12").Value = 320.75
